$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.796.67'
$ws.Range("E2").Value = '  +2.10%  '
$ws.Range("D3").Value = '3.004.58'
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '511.74'
$ws.Range("E5").Value = '  +4.32%  '
$ws.Range("D6").Value = '139.28'
$ws.Range("E6").Value = '  +5.21%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.435'
$ws.Range("E8").Value = '  +4.00%  '
$ws.Range("D9").Value = '7.48'
$ws.Range("E9").Value = '  +4.96%  '
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").Value = '  +6.81%  '
$ws.Range("D11").Value = '0.358'
$ws.Range("E11").Value = '  +2.35%  '
$ws.Range("D12").Value = '0.130'
$ws.Range("E12").Value = '  +2.30%  '
$ws.Range("D13").Value = '3.517.04'
$ws.Range("E13").Value = '  +1.12%  '
$ws.Range("D14").Value = '25.97'
$ws.Range("E14").Value = '  +4.38%  '
$ws.Range("D15").Value = '0.0000156'
$ws.Range("E15").Value = '  +11.48%  '
$ws.Range("D16").Value = '56.830.20'
$ws.Range("E16").Value = '  +2.44%  '
$ws.Range("D17").Value = '3.002.22'
$ws.Range("E17").Value = '  +1.61%  '
$ws.Range("D18").Value = '5.95'
$ws.Range("E18").Value = '  +5.88%  '
$ws.Range("D19").Value = '12.56'
$ws.Range("E19").Value = '  +3.51%  '
$ws.Range("D20").Value = '7.86'
$ws.Range("E20").Value = '  +4.88%  '
$ws.Range("D21").Value = '327.62'
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").Value = '0.488'
$ws.Range("E23").Value = '  +4.98%  '
$ws.Range("D24").Value = '63.46'
$ws.Range("E24").Value = '  +5.04%  '
$ws.Range("D25").Value = '0.171'
$ws.Range("E25").Value = '  +4.37%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("D27").Value = '0.0₃0912'
$ws.Range("E27").Value = '  +7.55%  '
$ws.Range("D28").Value = '6.67'
$ws.Range("E28").Value = '  +1.84%  '
$ws.Range("D29").Value = '7.05'
$ws.Range("E29").Value = '  +7.62%  '
$ws.Range("D30").Value = '1.23'
$ws.Range("E30").Value = '  +5.10%  '
$ws.Range("D31").Value = '1.81'
$ws.Range("E31").Value = '  +6.66%  '
$ws.Range("D32").Value = '20.71'
$ws.Range("E32").Value = '  +6.36%  '
$ws.Range("D33").Value = '154.75'
$ws.Range("E33").Value = '  +3.66%  '
$ws.Range("D34").Value = '4.57'
$ws.Range("E34").Value = '  +3.73%  '
$ws.Range("D35").Value = '5.72'
$ws.Range("E35").Value = '  +0.42%  '
$ws.Range("D36").Value = '1.28'
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("D37").Value = '0.0681'
$ws.Range("E37").Value = '  +4.70%  '
$ws.Range("D38").Value = '23.95'
$ws.Range("E38").Value = '  +2.58%  '
$ws.Range("D39").Value = '3.035.22'
$ws.Range("E39").Value = '  +1.19%  '
$ws.Range("D40").Value = '37.06'
$ws.Range("E40").Value = '  +2.54%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("D42").Value = '2.291.94'
$ws.Range("E42").Value = '  +7.53%  '
$ws.Range("D43").Value = '0.647'
$ws.Range("E43").Value = '  +2.56%  '
$ws.Range("D44").Value = '3.69'
$ws.Range("E44").Value = '  +4.07%  '
$ws.Range("D45").Value = '1.01'
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("D46").Value = '1.42'
$ws.Range("E46").Value = '  +2.56%  '
$ws.Range("D47").Value = '1.96'
$ws.Range("E47").Value = '  +11.43%  '
$ws.Range("D48").Value = '5.89'
$ws.Range("E48").Value = '  +5.39%  '
$ws.Range("D49").Value = '0.0239'
$ws.Range("E49").Value = '  +1.79%  '
$ws.Range("D50").Value = '19.35'
$ws.Range("E50").Value = '  +0.58%  '
$ws.Range("D51").Value = '0.0873'
$ws.Range("E51").Value = '  +4.99%  '
